$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the previously-last row (row 75, "01-01-2021") with revised figures ---
$ws.Range("B75").Value = 2.4
$ws.Range("C75").Value = 2.3
$ws.Range("D75").Value = 109.3
$ws.Range("E75").Value = 152.9

# --- Append the new quarter row (row 76) for "01-04-2021" ---
# Excel's COM layer auto-detects date-like strings (e.g. "01-04-2021") and
# silently stores them as date serials when assigned straight to .Value.
# Build the label as a text formula in an unused scratch cell (which stores
# it as a real string), then copy/paste-special the evaluated value into the
# target cell so it lands as a plain shared string, matching how the other
# period labels in column A are stored.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""01-04-2021"""
$scratch.Copy()
$ws.Range("A76").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("B76").Value = 2.7
$ws.Range("C76").Value = 2.6
$ws.Range("D76").Value = 107.6
$ws.Range("E76").Value = 155.2
